# Clear the GARANTI (column C) benchmark values for the EFT/HAVALE/SWIFT
# rows — the figures were stale/incorrect and are being removed pending
# re-verification. Style formatting (fill/border/number format) on the
# cells is left intact; only the text content is cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3,4,5,6,8,9,10,11,12,13,14)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).ClearContents()
}
